# Generate Report for Archive
#
# 1) The handback status text moves from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F4 and the
#    "Status" column on the per-locale sheets).
# 2) Because the new text is shorter, the Status-ish columns are
#    re-sized (narrower) on all three sheets.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value()
        if (("" + $v) -eq $oldText) {
            $cell.Value = $newText
        }
    }
}

# Narrow the columns that held the long status text to their new
# (shorter-content) width.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Columns.Item(5).ColumnWidth = 12.5   # column E
$ws1.Columns.Item(6).ColumnWidth = 12.5   # column F

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Columns.Item(3).ColumnWidth = 12.5   # column C (Status)

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Columns.Item(3).ColumnWidth = 12.5   # column C (Status)
